$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = '5840917 - Fabrício Maciel Gomes'
$ws.Range("C10").Value = '5840917 - Fabrício Maciel Gomes'

$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60.0

$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Introduction to Operational Research, Linear Programming, Simplex Method, Introduction to Graphs and Network Optimization, Case Study in Linear Programming, Introduction to Queue Theory.'
$ws.Range("C14").Value = 'Introduction to Operational Research, Linear Programming, Simplex Method, Introduction to Graphs and Network Optimization, Case Study in Linear Programming, Introduction to Queue Theory.'

$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2021'
$ws.Range("C15").Value = '01/01/2021'
$ws.Rows.Item(15).RowHeight = 120.0

$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models'
$ws.Range("C16").Value = '1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models'

$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).RowHeight = $ws.StandardHeight

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5840917 - Fabrício Maciel Gomes'
$ws.Range("C18").Value = '5840917 - Fabrício Maciel Gomes'
$ws.Rows.Item(18).RowHeight = 60.0

$ws.Range("A19").Value = 'Critério:'

$ws.Range("A20").Value = 'Norma de recuperação:'

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Rows.Item(21).RowHeight = 120.0

$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).RowHeight = $ws.StandardHeight

$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = 'LOB1012 -  Estatística  (Requisito fraco)`n'
$ws.Range("C23").Value = 'LOB1012 -  Estatística  (Requisito fraco)`n'
$ws.Rows.Item(23).RowHeight = 30.0

# Remove the now-obsolete last row (content shifted up by one)
$ws.Rows.Item(24).Delete()
